# Update the "Overview" sheet: roll the reporting years forward by one
# (drop 1396/12, shift 1397..1400 left, add new 1401/12 column) and refresh
# the underlying yearly figures ("update database and change read_price
# algorithm").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# --- Year header labels (row 8 and row 24) -------------------------------
# Columns E:I previously read 1396/12, 1397/12, 1398/12, 1399/12, 1400/12.
# They now read 1397/12, 1398/12, 1399/12, 1400/12, 1401/12.
$years = @(
    "دوازده ماهه منتهی به 1397/12",
    "دوازده ماهه منتهی به 1398/12",
    "دوازده ماهه منتهی به 1399/12",
    "دوازده ماهه منتهی به 1400/12",
    "دوازده ماهه منتهی به 1401/12"
)
$cols = @("E", "F", "G", "H", "I")

for ($i = 0; $i -lt 5; $i++) {
    $ws.Range($cols[$i] + "8").Value = $years[$i]
    $ws.Range($cols[$i] + "24").Value = $years[$i]
}

# --- Refreshed yearly data -------------------------------------------------
# row => values for columns E, F, G, H, I
$data = @{
    10 = @(0, 0, 0, 538109, 1437718)
    13 = @(4376, 2889, 1713, 1292, 3148)
    15 = @(536, 281, 170, 696, 221)
    16 = @(5494, 2656, 2198, 3013, 9119)
    17 = @(40332, 54575, 61561, 83724, 165037)
    19 = @(50368, 263208, 822655, 106160, 191655)
    20 = @(101106, 323609, 888297, 732994, 1806898)
    26 = @(422, 409, 410, 557, 580)
    27 = @(87, 86, 66, 144, 157)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt 5; $i++) {
        $ws.Range($cols[$i] + $row).Value = $vals[$i]
    }
}

# --- Theme accent color swap (accent1 <-> accent5) ------------------------
$theme = $wb.Theme
$colorScheme = $theme.ThemeColorScheme
$accent1 = $colorScheme.Colors(5).RGB
$accent5 = $colorScheme.Colors(9).RGB
$colorScheme.Colors(5).RGB = $accent5
$colorScheme.Colors(9).RGB = $accent1
